$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing table content
$ws.Range("A1:C15").ClearContents()

# New header + data (sorted by locus, A2:A9 ascending)
$data = @(
    @("locus", "Notes", "Keep or Drop"),
    @("Mnov_gtseq_110", "This is generally messy. It's amplifying two different products - the target (90bp) and a 127 bp off-target read", ""),
    @("Mnov_gtseq_155", "This locus probably has 3 good SNPs, but microhaplot is calling 5; might be recoverable", "??"),
    @("Mnov_gtseq_219", "This should be recoverable - lots of Ns due to masked reads", "??"),
    @("Mnov_gtseq_302", "There's a problem with the reverse primer; maybe can use forward reads only, but there might be an analogous locus being amplified", "??"),
    @("Mnov_gtseq_420", "This has almost no reads for a well-performing sample; check in rest of samples", ""),
    @("Mnov_gtseq_474", "This is messy and the vcf has waaaaayyyyy too many variants; might be recoverable if I figure out which sites are reliable SNPs", ""),
    @("Mnov_gtseq_535", "The reverse primer produced 28,188 reads of something not the target!!!", ""),
    @("Mnov_gtseq_541", "This is amplifying two different products in non-humpbacks", "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    if ($data[$i][2] -ne "") {
        $ws.Cells.Item($row, 3).Value = $data[$i][2]
    }
}

$sortRange = $ws.Range("A2:C9")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A9"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = $false
$ws.Sort.Apply()

$ws.Range("B14").Select() | Out-Null
